$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.533.19"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "1.913.34"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5237"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3963"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09661"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.534"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "1.916.13"
$ws.Range("E14").Value = "  +5.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.564"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06652"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.61%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.334"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.43%  "
$ws.Range("D23").Value = "28.649.61"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.693"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.13%  "
$ws.Range("D27").Value = "2.134.46"
$ws.Range("E27").Value = "  +5.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "159.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.106"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1087"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.752"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.642"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.892"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06747"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02438"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.257"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.46%  "
$ws.Range("E39").Value = "  +4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.103"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6445"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.188"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6091"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.769"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.284"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.035"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.211"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.41%  "
